# Edit sheet Card24 by admin
# Append one new blank row (row 13) to the bottom of the service-log table
# on the "Card24" sheet, extending the used range from A1:O12 to A1:O13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# The table currently spans A1:O12 (header row + 11 data rows). Touching the
# style of the next row (A13:O13) materializes those 15 cells as blank
# entries and grows the sheet's used range/dimension to A1:O13, matching a
# newly appended (still-empty) row in the log.
$ws.Range("A13:O13").Style = "Normal"
